$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").Value = 0.06458317429431877
$ws.Range("Y2").Value = 0.04706731114257444
$ws.Range("AB2").Value = 0.064577048110998
$ws.Range("AC2").Value = -0.06361135560782138
$ws.Range("X3").Value = 0.07798968720672064
$ws.Range("Y3").Value = 0.6690965598795264
$ws.Range("AB3").Value = 0.07605482288097691
$ws.Range("AC3").Value = 1.05945605376733
$ws.Range("X4").Value = 0.06462053199068069
$ws.Range("Y4").Value = 0.6482507551380323
$ws.Range("AB4").Value = 0.06460869610473566
$ws.Range("AC4").Value = 0.7538351079298465
$ws.Range("X5").Value = 0.06454268619014099
$ws.Range("Y5").Value = 0.1150529020451531
$ws.Range("AB5").Value = 0.06454268619014099
$ws.Range("AC5").Value = 0.3775387617736599
$ws.Range("X6").Value = 0.06454268619014099
$ws.Range("Y6").Value = 0.1077095660621113
$ws.Range("AB6").Value = 0.06454268619014099
$ws.Range("AC6").Value = 0.3585007920707286
$ws.Range("X7").Value = 0.06454268619014099
$ws.Range("Y7").Value = 0.07993891916102956
$ws.Range("AB7").Value = 0.06454268619014099
$ws.Range("AC7").Value = 0.3046880830406282
$ws.Range("X8").Value = 0.06454268619014099
$ws.Range("Y8").Value = 0.07461588986163896
$ws.Range("AB8").Value = 0.06454268619014099
$ws.Range("AC8").Value = 0.3008964072942784
$ws.Range("X9").Value = 0.06454268619014099
$ws.Range("Y9").Value = 0.07091548114053629
$ws.Range("AB9").Value = 0.06454268619014099
$ws.Range("AC9").Value = 0.2135321801200195
$ws.Range("X10").Value = 0.06454268619014099
$ws.Range("Y10").Value = 0.0393483643935166
$ws.Range("AB10").Value = 0.06454268619014099
$ws.Range("AC10").Value = 0.1408419291944744
$ws.Range("X11").Value = 0.06459534459298241
$ws.Range("Y11").Value = 0.021332160737508
$ws.Range("AB11").Value = 0.06458736433611648
$ws.Range("AC11").Value = 0.06101739934084668
$ws.Range("X12").Value = 0.06496410896764865
$ws.Range("Y12").Value = 0.005785891032351345
$ws.Range("AB12").Value = 0.06492879153429576
$ws.Range("AC12").Value = 0.04841862020492885
$ws.Range("X13").Value = 0.06497886617251346
$ws.Range("Y13").Value = -0.04673275584719521
$ws.Range("AB13").Value = 0.06490950320903179
$ws.Range("AC13").Value = -0.02681666768935469
$ws.Range("X14").Value = 0.06572863009550732
$ws.Range("Y14").Value = 0.1080994949044927
$ws.Range("AB14").Value = 0.06552314758884827
$ws.Range("AC14").Value = -0.03786444316320532
$ws.Range("X15").Value = 0.06458317429431877
$ws.Range("Y15").Value = -0.0633305203665056
$ws.Range("AB15").Value = 0.064577048110998
$ws.Range("AC15").Value = -0.06361135560782138
$ws.Range("X16").Value = 0.06517313416143791
$ws.Range("Y16").Value = 0.7127638285892783
$ws.Range("AB16").Value = 0.06507052208078896
$ws.Range("AC16").Value = -0.06507052208078896
$ws.Range("X17").Value = 0.06454268619014099
$ws.Range("Y17").Value = 0.1287517744512585
$ws.Range("AB17").Value = 0.06454268619014099
$ws.Range("AC17").Value = -0.0659565530557132
$ws.Range("X18").Value = 0.08848811219255928
$ws.Range("Y18").Value = -0.04142928866314752
$ws.Range("AB18").Value = 0.07761793699187824
$ws.Range("AC18").Value = -0.06918713212795236
$ws.Range("X19").Value = 0.07698850103744342
$ws.Range("Y19").Value = 0.6143087075011444
$ws.Range("AB19").Value = 0.07347655882122998
$ws.Range("AC19").Value = -0.07347655882122998
$ws.Range("X20").Value = 0.09131075211635725
$ws.Range("Y20").Value = -0.06165671587253517
$ws.Range("AB20").Value = 0.07856862968220912
$ws.Range("AC20").Value = -0.07734410572690865
$ws.Range("X21").Value = 0.06454268619014099
$ws.Range("Y21").Value = 0.3239708273233725
$ws.Range("AB21").Value = 0.06454268619014099
$ws.Range("AC21").Value = -0.07896111799987715
$ws.Range("X22").Value = 0.06454268619014099
$ws.Range("Y22").Value = 0.02299625461982786
$ws.Range("AB22").Value = 0.06454268619014099
$ws.Range("AC22").Value = -0.0845180042817484
$ws.Range("X23").Value = 0.06454987972857643
$ws.Range("Y23").Value = 0.04710060570831678
$ws.Range("AB23").Value = 0.0645507652073699
$ws.Range("AC23").Value = -0.08780657916085827
$ws.Range("X24").Value = 0.06465712388478528
$ws.Range("Y24").Value = -0.06507034702528115
$ws.Range("AB24").Value = 0.06464827840154172
$ws.Range("AC24").Value = -0.09979721457175449
$ws.Range("X25").Value = 0.06455734241300776
$ws.Range("Y25").Value = 0.008751390674323109
$ws.Range("AB25").Value = 0.0645591436509429
$ws.Range("AC25").Value = -0.1004591188923393
$ws.Range("X26").Value = 0.06454703864420445
$ws.Range("Y26").Value = -0.08106080011209436
$ws.Range("AB26").Value = 0.06454671247962875
$ws.Range("AC26").Value = -0.1889706756132692
$ws.Range("X27").Value = 0.066589565643368
$ws.Range("Y27").Value = -0.5353046494422506
$ws.Range("AB27").Value = 0.0667376399259958
$ws.Range("AC27").Value = -0.9587890109777139
